$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GLOBAL RESULTS")

# Insert the 4 new rows, bottom-to-top so earlier row numbers stay valid
$ws.Rows.Item(33).Insert()
$ws.Rows.Item(25).Insert()
$ws.Rows.Item(19).Insert()
$ws.Rows.Item(10).Insert()

# Populate the newly inserted rows
$ws.Range("A10").Value = "Actual Passengers Mass"
$ws.Range("B10").Value = "kg"
$ws.Range("C10").Value = 6460.0

$ws.Range("A20").Value = "Furnishings and Equipments Mass"
$ws.Range("B20").Value = "kg"
$ws.Range("C20").Value = 1251.8725753832337

$ws.Range("A27").Value = "Actual Passengers Weight"
$ws.Range("B27").Value = "N"
$ws.Range("C27").Value = 63350.95899999998

$ws.Range("A36").Value = "Furnishings and Equipments Weight"
$ws.Range("B36").Value = "N"
$ws.Range("C36").Value = 12276.676191381986

# Update recalculated totals downstream of the new rows
$ws.Range("C6").Value = 21955.92750895614  # Maximum Take-Off Mass
$ws.Range("C7").Value = 21575.92750895614  # Take-Off Mass
$ws.Range("C8").Value = 21297.249683687452  # Maximum Landing Mass
$ws.Range("C12").Value = 3036.3083506198227  # Fuel Mass
$ws.Range("C14").Value = 18919.619158336318  # Maximum Zero Fuel Mass
$ws.Range("C15").Value = 18539.619158336318  # Zero Fuel Mass
$ws.Range("C16").Value = 12079.619158336322  # Operating Empty Mass
$ws.Range("C17").Value = 11850.07551283632  # Empty Mass
$ws.Range("C18").Value = 11229.65151283632  # Manufacturer Empty Mass
$ws.Range("C23").Value = 215314.09650570469  # Maximum Take-Off Weight
$ws.Range("C24").Value = 211587.56950570468  # Take-Off Weight
$ws.Range("C25").Value = 208854.67361053347  # Maximum Landing Weight
$ws.Range("C30").Value = 185538.0832190988  # Maximum Zero Fuel Weight
$ws.Range("C31").Value = 181811.5562190988  # Zero Fuel Weight
$ws.Range("C32").Value = 118460.59721909885  # Operating Empty Weight
$ws.Range("C33").Value = 116209.54302795627  # Empty Weight
$ws.Range("C34").Value = 110125.26200835628  # Manufacturer Empty Weight
